# Add a new "Reg Proc" column to the requirements table on the "Details"
# sheet, and populate clarification/research notes that go with it (plus an
# update to an existing "Comments" cell for row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")
$tbl = $ws.ListObjects.Item(1)

# 1. Grow the table by one column (A2:S16 -> A2:T16) and name the header
#    cell to match the rest of the header row's formatting.
$newCol = $tbl.ListColumns.Add()
$ws.Range("T2").Value = "Reg Proc"
$ws.Range("S2").Copy()
$ws.Range("T2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. Fill in the new column's notes, in the same order the original author
#    typed them (first use of each string fixes its shared-string slot).
$ws.Range("T5").Value = "No Mapping of such kind from Reg Processor`nID Repo- Might not be there in ID Repo as well"
$ws.Range("T6").Value = "ID Repo- need to know "
$ws.Range("T9").Value = "Under processing`nProcessed"
$ws.Range("T8").Value = "Under processing`nProcessed`n"
$ws.Range("S8").Value = "Reg proc`nArchival policy"
$ws.Range("T10").Value = "E-UIN Generation"
$ws.Range("T7").Value = "there shud be a label as Res_Service`nReg Client packet needs to be understood`nService from Reg proc needs to be developed"
$ws.Range("T4").Value = "When UIN IS needed to be generated`n1.the Acknowledgment from Print queue- what needs to be done`nTime period `n2. If there is a print failure- no need to handle from MOSIP`nUser Story ?"

# 3. Match cell formatting: most new cells wrap their text like the
#    neighbouring "Module Dependency"/"Research info" columns; a couple stay
#    unwrapped (T6, T10). S8 switches from unwrapped to wrapped as well.
foreach ($addr in @("T4", "T5", "T7", "T8", "T9", "S8")) {
    $cell = $ws.Range($addr)
    $cell.WrapText = $true
    $cell.HorizontalAlignment = -4131
    $cell.VerticalAlignment = -4160
}

# 4. New column width, sized to fit its content like the author left it.
$ws.Range("T1").EntireColumn.ColumnWidth = 32.08984375

# 5. Leave the view the way the workbook was last saved: scrolled so the
#    frozen pane's corner shows N3, with T4 selected as the active cell.
$ws.Range("T4").Select()
$excel.ActiveWindow.ScrollColumn = 14
$excel.ActiveWindow.ScrollRow = 3
